# "Add files via upload" - append Customer/Distributor/Retailer/Employee
# Q&A rows to the chatbot intents sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("A76").Value = "Customer"
$ws.Range("B76").Value = "Great! What product are you looking for today?"

$ws.Range("A77").Value = "Distributor"
$ws.Range("B77").Value = "Great! Please enter your Crompton UserID"

$ws.Range("A78").Value = "Retailer"
$ws.Range("B78").Value = "Great! Please enter your Crompton UserID"

$ws.Range("A79").Value = "Employee"
$ws.Range("B79").Value = "Great! Please enter your Crompton UserID"

# Match the author's saved view: scrolled down with B79 selected.
$ws.Range("B79").Select()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
